# Insert a new row for the "uploaded-date" field into the "documents"
# sub-object of the Application component (becomes new row 12), pushing
# every subsequent row down by one. Excel's native Rows.Insert() keeps all
# existing cell values/styles/merges correctly shifted; we only need to
# populate the freshly inserted row and extend the one merged-cell pair
# whose range boundary lands exactly on the insertion point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 12 (shifts rows 12:82 down to 13:83).
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the "uploaded-date" field, mirroring the
# surrounding "documents" rows for the repeated/shared columns C:G.
$ws.Range("C12").Value = "The details of the application payload to be submitted"
$ws.Range("D12").Value = "application"
$ws.Range("E12").Value = "Application"
$ws.Range("F12").Value = "documents"
$ws.Range("G12").Value = "Documents[]"
$ws.Range("H12").Value = "uploaded-date"
$ws.Range("I12").Value = "Uploaded date"
$ws.Range("L12").Value = "The date the document was uploaded to the application"
$ws.Range("M12").Value = "date"
$ws.Range("N12").Value = "MUST"

# The "Checklist" section header merge (single-row, A/B columns) sat at
# row 44 before the insert and should now cover row 45.
$ws.Range("A45").Merge()
$ws.Range("B45").Merge()
